# Insert a new row at row 262 on the active worksheet, shifting existing
# rows 262:386 down to 263:387, and populate the new row with the new
# weekly price record (Lane Late / Primera, date serial 44553).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 262.
$ws.Rows.Item(262).Insert()

# The constant columns (market / region / product classification / unit /
# origin / kg-per-unit) are identical for every record in this table, so
# copy them from the row directly below (old row 262, now shifted to 263).
$ws.Cells.Item(262, 1).Value2  = $ws.Cells.Item(263, 1).Value2   # A Mercado ID
$ws.Cells.Item(262, 2).Value2  = $ws.Cells.Item(263, 2).Value2   # B Mercado
$ws.Cells.Item(262, 3).Value2  = $ws.Cells.Item(263, 3).Value2   # C Region
$ws.Cells.Item(262, 4).Value2  = 44553                           # D Fecha
$ws.Cells.Item(262, 5).Value2  = $ws.Cells.Item(263, 5).Value2   # E Codreg
$ws.Cells.Item(262, 6).Value2  = $ws.Cells.Item(263, 6).Value2   # F Tipo
$ws.Cells.Item(262, 7).Value2  = $ws.Cells.Item(263, 7).Value2   # G Producto ID
$ws.Cells.Item(262, 8).Value2  = $ws.Cells.Item(263, 8).Value2   # H Producto
$ws.Cells.Item(262, 9).Value2  = $ws.Cells.Item(263, 9).Value2   # I Categoria ID
$ws.Cells.Item(262, 10).Value2 = $ws.Cells.Item(263, 10).Value2  # J Categoria
$ws.Cells.Item(262, 11).Value2 = "Lane Late"                     # K Variedad
$ws.Cells.Item(262, 12).Value2 = "Primera"                       # L Calidad
$ws.Cells.Item(262, 13).Value2 = 24                               # M Volumen
$ws.Cells.Item(262, 14).Value2 = 210000                           # N Precio minimo
$ws.Cells.Item(262, 15).Value2 = 220000                           # O Precio maximo
$ws.Cells.Item(262, 16).Value2 = 215000                           # P Precio promedio ponderado
$ws.Cells.Item(262, 17).Value2 = $ws.Cells.Item(263, 17).Value2   # Q Unidad de comercializacion
$ws.Cells.Item(262, 18).Value2 = $ws.Cells.Item(263, 18).Value2   # R Origen
$ws.Cells.Item(262, 19).Value2 = 538                               # S Precio $/Kg
$ws.Cells.Item(262, 20).Value2 = $ws.Cells.Item(263, 20).Value2   # T Kg / unidad
